$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): want-to-go counts updated
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1893
$ws1.Range("F4").Value = 824
$ws1.Range("F5").Value = 785
$ws1.Range("F6").Value = 243

# Sheet "全部类型" (sheet4): same counts updated for matching events
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1893
$ws4.Range("F5").Value = 824
$ws4.Range("F6").Value = 785
$ws4.Range("F7").Value = 243
